$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1310.6
$ws.Range("J17").Value = 1341.1228
$ws.Range("L17").Value = 4023.3684
$ws.Range("N17").Value = -4359.3684
$ws.Range("H19").Value = 1038.8667
$ws.Range("I19").Value = 389.75
$ws.Range("K19").Value = 389.75
$ws.Range("M19").Value = -214.75
$ws.Range("H87").Value = 66428.57000000001
$ws.Range("J87").Value = 70000
$ws.Range("L87").Value = 70000
$ws.Range("N87").Value = -72496
$ws.Range("H90").Value = 66428.57000000001
$ws.Range("J90").Value = 70000
$ws.Range("L90").Value = 210000
$ws.Range("N90").Value = -222480
$ws.Range("H112").Value = 1258.7903
$ws.Range("J112").Value = 1258.7903
$ws.Range("L112").Value = 3776.3709
$ws.Range("N112").Value = -5992.3709
$ws.Range("H116").Value = 2847.4
$ws.Range("I116").Value = 3395.2
$ws.Range("J116").Value = 2299.6
$ws.Range("K116").Value = 3395.2
$ws.Range("L116").Value = 2299.6
$ws.Range("M116").Value = 46.80000000000018
$ws.Range("N116").Value = -9183.6
$ws.Range("H132").Value = 2463.4119
$ws.Range("I132").Value = 2250.1936
$ws.Range("K132").Value = 6750.5808
$ws.Range("M132").Value = -4220.5808
$ws.Range("H135").Value = 1539.9166
$ws.Range("I135").Value = 1539.9166
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 13859.2494
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -11324.2494
$ws.Range("N135").ClearContents()
$ws.Range("H138").Value = 13159987
$ws.Range("I138").Value = 41667904
$ws.Range("J138").Value = 2485.6345
$ws.Range("K138").Value = 125003712
$ws.Range("L138").Value = 7456.9035
$ws.Range("M138").Value = -124998572
$ws.Range("N138").Value = -17736.9035

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2431.8984
$ws.Range("I32").Value = 1535.2
$ws.Range("K32").Value = 1535.2
$ws.Range("M32").Value = -1248.2
$ws.Range("H61").Value = 3749.6667
$ws.Range("I61").Value = 1750
$ws.Range("J61").Value = 4749.5
$ws.Range("K61").Value = 1750
$ws.Range("L61").Value = 4749.5
$ws.Range("M61").Value = -1538
$ws.Range("N61").Value = -5173.5
$ws.Range("H74").Value = 2623.4912
$ws.Range("I74").Value = 2165.0715
$ws.Range("K74").Value = 2165.0715
$ws.Range("M74").Value = -1291.0715
$ws.Range("H77").Value = 2623.4912
$ws.Range("I77").Value = 2165.0715
$ws.Range("K77").Value = 10825.3575
$ws.Range("M77").Value = -6457.3575
$ws.Range("H103").Value = 40362
$ws.Range("J103").Value = 40362
$ws.Range("L103").Value = 40362
$ws.Range("N103").Value = -42706
$ws.Range("H132").Value = 7339.686
$ws.Range("I132").Value = 5106.8716
$ws.Range("J132").Value = 17108.25
$ws.Range("K132").Value = 15320.6148
$ws.Range("L132").Value = 51324.75
$ws.Range("M132").Value = -12790.6148
$ws.Range("N132").Value = -56384.75
$ws.Range("H136").Value = 3749.6667
$ws.Range("I136").Value = 1750
$ws.Range("J136").Value = 4749.5
$ws.Range("K136").Value = 5250
$ws.Range("L136").Value = 14248.5
$ws.Range("M136").Value = -2700
$ws.Range("N136").Value = -19348.5

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 15299.8
$ws.Range("I5").Value = 15299.8
$ws.Range("K5").Value = 15299.8
$ws.Range("M5").Value = -15186.8
$ws.Range("H20").Value = 2381
$ws.Range("I20").Value = 2359.5
$ws.Range("J20").Value = 2415.4
$ws.Range("K20").Value = 2359.5
$ws.Range("L20").Value = 2415.4
$ws.Range("M20").Value = -2112.5
$ws.Range("N20").Value = -2909.4
$ws.Range("H99").Value = 200001470
$ws.Range("I99").Value = 250001120
$ws.Range("K99").Value = 250001120
$ws.Range("M99").Value = -249999622
$ws.Range("H134").Value = 3679.5356
$ws.Range("I134").Value = 2815.75
$ws.Range("J134").Value = 5839
$ws.Range("K134").Value = 8447.25
$ws.Range("L134").Value = 17517
$ws.Range("M134").Value = -5912.25
$ws.Range("N134").Value = -22587

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 711
$ws.Range("I10").Value = 300
$ws.Range("J10").Value = 1122
$ws.Range("K10").Value = 300
$ws.Range("L10").Value = 1122
$ws.Range("M10").Value = -161
$ws.Range("N10").Value = -1400
$ws.Range("H31").Value = 1613.2727
$ws.Range("I31").Value = 1404.6666
$ws.Range("K31").Value = 1404.6666
$ws.Range("M31").Value = -1109.6666
$ws.Range("H34").Value = 1613.2727
$ws.Range("I34").Value = 1404.6666
$ws.Range("K34").Value = 1404.6666
$ws.Range("M34").Value = -1202.6666
$ws.Range("H58").Value = 3670.75
$ws.Range("I58").Value = 4061.0667
$ws.Range("K58").Value = 4061.0667
$ws.Range("M58").Value = -3858.0667
$ws.Range("H99").Value = 4445.154
$ws.Range("I99").Value = 1942
$ws.Range("J99").Value = 7365.5
$ws.Range("K99").Value = 1942
$ws.Range("L99").Value = 7365.5
$ws.Range("M99").Value = -444
$ws.Range("N99").Value = -10361.5
$ws.Range("H122").Value = 2749.923
$ws.Range("I122").Value = 2470.875
$ws.Range("K122").Value = 7412.625
$ws.Range("M122").Value = -4962.625
$ws.Range("H126").Value = 4445.154
$ws.Range("I126").Value = 1942
$ws.Range("J126").Value = 7365.5
$ws.Range("K126").Value = 5826
$ws.Range("L126").Value = 22096.5
$ws.Range("M126").Value = -3356
$ws.Range("N126").Value = -27036.5
$ws.Range("H132").Value = 1486.6364
$ws.Range("I132").Value = 1486.6364
$ws.Range("K132").Value = 4459.9092
$ws.Range("M132").Value = -1929.9092
$ws.Range("H134").Value = 1623.7646
$ws.Range("I134").Value = 1623.7646
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4871.293799999999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2336.293799999999
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 3670.75
$ws.Range("I136").Value = 4061.0667
$ws.Range("K136").Value = 12183.2001
$ws.Range("M136").Value = -9633.2001

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 108.4
$ws.Range("I12").Value = 184.5
$ws.Range("J12").Value = 57.666668
$ws.Range("K12").Value = 553.5
$ws.Range("L12").Value = 173.000004
$ws.Range("M12").Value = -380.5
$ws.Range("N12").Value = -519.000004
$ws.Range("H63").Value = 10156
$ws.Range("J63").Value = 8564.25
$ws.Range("L63").Value = 25692.75
$ws.Range("N63").Value = -27190.75
$ws.Range("H64").Value = 11837.667
$ws.Range("I64").Value = 7750
$ws.Range("J64").Value = 20013
$ws.Range("K64").Value = 23250
$ws.Range("L64").Value = 60039
$ws.Range("M64").Value = -22980
$ws.Range("N64").Value = -60579
$ws.Range("H66").Value = 10156
$ws.Range("J66").Value = 8564.25
$ws.Range("L66").Value = 77078.25
$ws.Range("N66").Value = -84566.25
$ws.Range("H67").Value = 11837.667
$ws.Range("I67").Value = 7750
$ws.Range("J67").Value = 20013
$ws.Range("K67").Value = 23250
$ws.Range("L67").Value = 60039
$ws.Range("M67").Value = -22314
$ws.Range("N67").Value = -61911
$ws.Range("H119").Value = 3566.6667
$ws.Range("I119").Value = 3566.6667
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 10700.0001
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = -5862.000100000001
$ws.Range("N119").ClearContents()

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2977.375
$ws.Range("I102").Value = 3017
$ws.Range("J102").Value = 2700
$ws.Range("K102").Value = 3017
$ws.Range("L102").Value = 2700
$ws.Range("M102").Value = -1395
$ws.Range("N102").Value = -5944
$ws.Range("H126").Value = 3485.2666
$ws.Range("I126").Value = 3044.5386
$ws.Range("J126").Value = 6350
$ws.Range("K126").Value = 9133.6158
$ws.Range("L126").Value = 19050
$ws.Range("M126").Value = -6663.6158
$ws.Range("N126").Value = -23990
$ws.Range("H132").Value = 2449.4644
$ws.Range("I132").Value = 2449.4644
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7348.3932
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4818.3932
$ws.Range("N132").ClearContents()

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6331.4116
$ws.Range("I122").Value = 5249.6665
$ws.Range("K122").Value = 15748.9995
$ws.Range("M122").Value = -13298.9995
$ws.Range("H132").Value = 3930.4634
$ws.Range("I132").Value = 3418.6667
$ws.Range("J132").Value = 4917.5
$ws.Range("K132").Value = 10256.0001
$ws.Range("L132").Value = 14752.5
$ws.Range("M132").Value = -7726.000100000001
$ws.Range("N132").Value = -19812.5

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H81").Value = 2021.5714
$ws.Range("I81").Value = 2021.5714
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 4043.1428
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -2982.1428
$ws.Range("N81").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H84").Value = 2021.5714
$ws.Range("I84").Value = 2021.5714
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 20215.714
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -14911.714
$ws.Range("N84").ClearContents()
$ws.Range("H113").Value = 1252.9143
$ws.Range("I113").Value = 762
$ws.Range("K113").Value = 2286
$ws.Range("M113").Value = -116
$ws.Range("H136").Value = 4503.5713
$ws.Range("I136").Value = 4718.3447
$ws.Range("K136").Value = 14155.0341
$ws.Range("M136").Value = -11605.0341
